$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections in rows above the deleted row (rows 3-25) ---

# Row 3: D3 was missing -> -14.2
$ws.Range("D3").Value = -14.2

# Row 4: E4 was -6.4 -> now missing
$ws.Range("E4").ClearContents()

# Row 5: D5 was -14.4 -> now missing
$ws.Range("D5").ClearContents()

# Row 6: F6 was missing -> 16.43
$ws.Range("F6").Value = 16.43

# Row 9: E9 was missing -> -6.8
$ws.Range("E9").Value = -6.8

# Row 10: E10 was missing -> -6.1
$ws.Range("E10").Value = -6.1

# Row 12: F12 was 17.45 -> now missing
$ws.Range("F12").ClearContents()

# Row 13: E13 was -5.3 -> now missing
$ws.Range("E13").ClearContents()

# Row 14: E14 was -5.4 -> now missing; F14 was missing -> 17.76
$ws.Range("E14").ClearContents()
$ws.Range("F14").Value = 17.76

# Row 19: F19 was missing -> 17.81
$ws.Range("F19").Value = 17.81

# Row 20: F20 was 17.73 -> now missing
$ws.Range("F20").ClearContents()

# Row 21: D21 was missing -> -14.3
$ws.Range("D21").Value = -14.3

# Row 23: D23 was -13.9 -> now missing
$ws.Range("D23").ClearContents()

# Row 25: F25 was 16.6 -> now missing
$ws.Range("F25").ClearContents()

# --- Remove row 26 (RM 232) and row 28 (SC 92) entirely, shifting the rest up ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Fill in previously-missing values for rows that shifted up ---
# Former row 29 (SC 101) is now row 27; its F value becomes 17
$ws.Range("F27").Value = 17

# Former row 30 (SC 105) is now row 28; its F value becomes 17.44
$ws.Range("F28").Value = 17.44

# Former row 34 (SC 193) is now row 32; its D value becomes -14.7
$ws.Range("D32").Value = -14.7
